$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'21.716.13"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "'1.538.55"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'290.00"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.3939"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3202"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "'43.47"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "'0.07210"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "'1.078"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'5.772"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "'18.49"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").Value = "'6.641"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'1.540.93"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'0.00001097"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("D18").Value = "'0.06618"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'84.28"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'6.156"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "'15.56"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").Value = "'10.84"
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("D24").Value = "'2.367"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "'21.725.14"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "'2.395"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").Value = "'151.14"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'18.55"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'4.864"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'1.700.76"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'117.77"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'6.094"
$ws.Range("E32").Value = "  +6.29%  "
$ws.Range("D33").Value = "'0.9715"
$ws.Range("E33").Value = "  -8.63%  "
$ws.Range("D34").Value = "'0.08104"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").Value = "'5.213"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "'8.518"
$ws.Range("E36").Value = "  -7.68%  "
$ws.Range("D37").Value = "'1.496"
$ws.Range("E37").Value = "  -8.13%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("D40").Value = "'11.30"
$ws.Range("E40").Value = "  +4.69%  "
$ws.Range("D41").Value = "'0.2051"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'1.185"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'0.5828"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "'13.17"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "'0.5597"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "'1.894"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'116.07"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").Value = "'0.06727"
$ws.Range("E51").Value = "  -1.95%  "
